$d = $word.ActiveDocument

# Mara's line gets two extra exclamation points:
# "Mara: Alright, alright! I change my mind..." -> "Mara: Alright, alright!!! I change my mind..."
$d.Content.Find.Execute(
    "Mara: Alright, alright! I change my mind",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Mara: Alright, alright!!! I change my mind",
    2
)
